# Auto-generated edit script applying scheduled-runner value updates
# to the Tiamat_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (ALC) - leve item G18
$ws.Range("H18").Value = 244.56667
$ws.Range("I18").Value = 175.44444
$ws.Range("J18").Value = 866.6667
$ws.Range("K18").Value = 175.44444
$ws.Range("L18").Value = 866.6667
$ws.Range("M18").Value = 108.55556
$ws.Range("N18").Value = -1434.6667

# Row 111 (ALC) - leve item G111
$ws.Range("H111").Value = 2076.375
$ws.Range("I111").Value = 1596.5
$ws.Range("J111").Value = 3516
$ws.Range("K111").Value = 4789.5
$ws.Range("L111").Value = 10548
$ws.Range("M111").Value = -1722.5
$ws.Range("N111").Value = -16682

# Row 118 (ALC) - leve item G118
$ws.Range("H118").Value = 807.125
$ws.Range("I118").Value = 516.1
$ws.Range("J118").Value = 1015
$ws.Range("K118").Value = 1548.3
$ws.Range("L118").Value = 3045
$ws.Range("M118").Value = 108.6999999999998
$ws.Range("N118").Value = -6359

# Row 137 (ALC) - leve item G137
$ws.Range("H137").Value = 5831.9565
$ws.Range("I137").Value = 988.5
$ws.Range("K137").Value = 2965.5
$ws.Range("M137").Value = -415.5

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (ARM) - leve item G63
$ws.Range("H63").Value = 2646.6155
$ws.Range("I63").Value = 2754.5454
$ws.Range("J63").Value = 2053
$ws.Range("K63").Value = 2754.5454
$ws.Range("L63").Value = 2053
$ws.Range("M63").Value = -2068.5454
$ws.Range("N63").Value = -3425

# Row 66 (ARM) - leve item G66
$ws.Range("H66").Value = 2646.6155
$ws.Range("I66").Value = 2754.5454
$ws.Range("J66").Value = 2053
$ws.Range("K66").Value = 13772.727
$ws.Range("L66").Value = 10265
$ws.Range("M66").Value = -10340.727
$ws.Range("N66").Value = -17129

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (BSM) - leve item G99
$ws.Range("H99").Value = 1453.6842
$ws.Range("I99").Value = 1414.4445
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 1414.4445
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = 83.55549999999994
$ws.Range("N99").Value = -4546

# Row 105 (BSM) - leve item G105
$ws.Range("H105").Value = 885588.9
$ws.Range("I105").Value = 1593140
$ws.Range("J105").Value = 1150
$ws.Range("K105").Value = 1593140
$ws.Range("L105").Value = 1150
$ws.Range("M105").Value = -1591393
$ws.Range("N105").Value = -4644

# Row 134 (BSM) - leve item G134
$ws.Range("H134").Value = 19629154
$ws.Range("I134").Value = 1635.279
$ws.Range("J134").Value = 125127064
$ws.Range("K134").Value = 4905.837
$ws.Range("L134").Value = 375381192
$ws.Range("M134").Value = -2370.837
$ws.Range("N134").Value = -375386262

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP) - leve item G7
$ws.Range("H7").Value = 1071.1
$ws.Range("I7").Value = 1487.1428
$ws.Range("J7").Value = 100.333336
$ws.Range("K7").Value = 1487.1428
$ws.Range("L7").Value = 100.333336
$ws.Range("M7").Value = -1374.1428
$ws.Range("N7").Value = -326.333336

# Row 16 (CRP) - leve item G16
$ws.Range("H16").Value = 851.5714
$ws.Range("I16").Value = 665.25
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 665.25
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -378.25
$ws.Range("N16").Value = -1674

# Row 22 (CRP) - leve item G22
$ws.Range("H22").Value = 831.45
$ws.Range("I22").Value = 1017.7143
$ws.Range("J22").Value = 396.83334
$ws.Range("K22").Value = 1017.7143
$ws.Range("L22").Value = 396.83334
$ws.Range("M22").Value = -667.7143
$ws.Range("N22").Value = -1096.83334

# Row 31 (CRP) - leve item G31
$ws.Range("H31").Value = 4308.262
$ws.Range("I31").Value = 1302.359
$ws.Range("J31").Value = 9636.909
$ws.Range("K31").Value = 1302.359
$ws.Range("L31").Value = 9636.909
$ws.Range("M31").Value = -1007.359
$ws.Range("N31").Value = -10226.909

# Row 34 (CRP) - leve item G34
$ws.Range("H34").Value = 4308.262
$ws.Range("I34").Value = 1302.359
$ws.Range("J34").Value = 9636.909
$ws.Range("K34").Value = 1302.359
$ws.Range("L34").Value = 9636.909
$ws.Range("M34").Value = -1100.359
$ws.Range("N34").Value = -10040.909

# Row 86 (CRP) - leve item G86
$ws.Range("H86").Value = 500050000
$ws.Range("I86").Value = 500050000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 500050000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -500048877
$ws.Range("N86").ClearContents()

# Row 89 (CRP) - leve item G89
$ws.Range("H89").Value = 500050000
$ws.Range("I89").Value = 500050000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2500250000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2500244384
$ws.Range("N89").ClearContents()

# Row 105 (CRP) - leve item G105
$ws.Range("H105").Value = 1525.5555
$ws.Range("I105").Value = 1155
$ws.Range("J105").Value = 1631.4286
$ws.Range("K105").Value = 1155
$ws.Range("L105").Value = 1631.4286
$ws.Range("M105").Value = 592
$ws.Range("N105").Value = -5125.4286

# Row 113 (CRP) - leve item G113
$ws.Range("H113").Value = 851.5714
$ws.Range("I113").Value = 665.25
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 665.25
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1504.75
$ws.Range("N113").Value = -5440

# Row 118 (CRP) - leve item G118
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

$ws = $wb.Worksheets.Item("CUL")
# Row 97 (CUL) - leve item G97
$ws.Range("H97").Value = 4184.533
$ws.Range("I97").Value = 1790
$ws.Range("J97").Value = 4783.1665
$ws.Range("K97").Value = 5370
$ws.Range("L97").Value = 14349.4995
$ws.Range("M97").Value = -4874
$ws.Range("N97").Value = -15341.4995

# Row 98 (CUL) - leve item G98
$ws.Range("H98").Value = 337.66666
$ws.Range("I98").Value = 327.8
$ws.Range("J98").Value = 350
$ws.Range("K98").Value = 983.4000000000001
$ws.Range("L98").Value = 1050
$ws.Range("M98").Value = 514.5999999999999
$ws.Range("N98").Value = -4046

# Row 107 (CUL) - leve item G107
$ws.Range("H107").Value = 595.34784
$ws.Range("I107").Value = 406.75
$ws.Range("J107").Value = 801.0909
$ws.Range("K107").Value = 1220.25
$ws.Range("L107").Value = 2403.2727
$ws.Range("M107").Value = 699.75
$ws.Range("N107").Value = -6243.2727

# Row 122 (CUL) - leve item G122
$ws.Range("H122").Value = 650.1579
$ws.Range("I122").Value = 495.07693
$ws.Range("J122").Value = 986.1667
$ws.Range("K122").Value = 4455.69237
$ws.Range("L122").Value = 8875.5003
$ws.Range("M122").Value = -2005.69237
$ws.Range("N122").Value = -13775.5003

# Row 131 (CUL) - leve item G131
$ws.Range("H131").Value = 162077.84
$ws.Range("I131").Value = 395.94736
$ws.Range("J131").Value = 233518.67
$ws.Range("K131").Value = 1187.84208
$ws.Range("L131").Value = 700556.01
$ws.Range("M131").Value = 3852.15792
$ws.Range("N131").Value = -710636.01

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM) - leve item G70
$ws.Range("H70").Value = 3725.3333
$ws.Range("I70").Value = 3348
$ws.Range("J70").Value = 4102.6665
$ws.Range("K70").Value = 3348
$ws.Range("L70").Value = 4102.6665
$ws.Range("M70").Value = -3078
$ws.Range("N70").Value = -4642.6665

# Row 73 (GSM) - leve item G73
$ws.Range("H73").Value = 3725.3333
$ws.Range("I73").Value = 3348
$ws.Range("J73").Value = 4102.6665
$ws.Range("K73").Value = 3348
$ws.Range("L73").Value = 4102.6665
$ws.Range("M73").Value = -2412
$ws.Range("N73").Value = -5974.6665

# Row 80 (GSM) - leve item G80
$ws.Range("H80").Value = 5193.4473
$ws.Range("I80").Value = 2849.9285
$ws.Range("J80").Value = 6560.5
$ws.Range("K80").Value = 2849.9285
$ws.Range("L80").Value = 6560.5
$ws.Range("M80").Value = -1851.9285
$ws.Range("N80").Value = -8556.5

# Row 83 (GSM) - leve item G83
$ws.Range("H83").Value = 5193.4473
$ws.Range("I83").Value = 2849.9285
$ws.Range("J83").Value = 6560.5
$ws.Range("K83").Value = 14249.6425
$ws.Range("L83").Value = 32802.5
$ws.Range("M83").Value = -9257.6425
$ws.Range("N83").Value = -42786.5

# Row 96 (GSM) - leve item G96
$ws.Range("H96").Value = 7995
$ws.Range("J96").Value = 7995
$ws.Range("L96").Value = 7995
$ws.Range("N96").Value = -13487

# Row 97 (GSM) - leve item G97
$ws.Range("H97").Value = 1184.6471
$ws.Range("I97").Value = 1174.2142
$ws.Range("J97").Value = 1233.3334
$ws.Range("K97").Value = 1174.2142
$ws.Range("L97").Value = 1233.3334
$ws.Range("M97").Value = -678.2141999999999
$ws.Range("N97").Value = -2225.3334

# Row 122 (GSM) - leve item G122
$ws.Range("H122").Value = 3310.9092
$ws.Range("I122").Value = 2773
$ws.Range("J122").Value = 4252.25
$ws.Range("K122").Value = 8319
$ws.Range("L122").Value = 12756.75
$ws.Range("M122").Value = -5869
$ws.Range("N122").Value = -17656.75

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW) - leve item G16
$ws.Range("H16").Value = 1252.579
$ws.Range("I16").Value = 771.2143
$ws.Range("J16").Value = 2600.4
$ws.Range("K16").Value = 771.2143
$ws.Range("L16").Value = 2600.4
$ws.Range("M16").Value = -601.2143
$ws.Range("N16").Value = -2940.4

# Row 40 (LTW) - leve item G40
$ws.Range("H40").Value = 35006.453
$ws.Range("I40").Value = 2080
$ws.Range("J40").Value = 58786.668
$ws.Range("K40").Value = 2080
$ws.Range("L40").Value = 58786.668
$ws.Range("M40").Value = -1944
$ws.Range("N40").Value = -59058.668

# Row 132 (LTW) - leve item G132
$ws.Range("H132").Value = 238007.69
$ws.Range("I132").Value = 64117.594
$ws.Range("J132").Value = 530875.2
$ws.Range("K132").Value = 192352.782
$ws.Range("L132").Value = 1592625.6
$ws.Range("M132").Value = -189822.782
$ws.Range("N132").Value = -1597685.6

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR) - leve item G62
$ws.Range("H62").Value = 5100
$ws.Range("I62").Value = 4762.5
$ws.Range("J62").Value = 5400
$ws.Range("K62").Value = 4762.5
$ws.Range("L62").Value = 5400
$ws.Range("M62").Value = -4138.5
$ws.Range("N62").Value = -6648

# Row 65 (WVR) - leve item G65
$ws.Range("H65").Value = 5100
$ws.Range("I65").Value = 4762.5
$ws.Range("J65").Value = 5400
$ws.Range("K65").Value = 23812.5
$ws.Range("L65").Value = 27000
$ws.Range("M65").Value = -20692.5
$ws.Range("N65").Value = -33240

# Row 95 (WVR) - leve item G95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 96 (WVR) - leve item G96
$ws.Range("H96").Value = 2738.7693
$ws.Range("I96").Value = 2222.6667
$ws.Range("J96").Value = 3900
$ws.Range("K96").Value = 2222.6667
$ws.Range("L96").Value = 3900
$ws.Range("M96").Value = -849.6667000000002
$ws.Range("N96").Value = -6646

# Row 126 (WVR) - leve item G126
$ws.Range("H126").Value = 689.8
$ws.Range("I126").Value = 611.625
$ws.Range("J126").Value = 1002.5
$ws.Range("K126").Value = 1834.875
$ws.Range("L126").Value = 3007.5
$ws.Range("M126").Value = 635.125
$ws.Range("N126").Value = -7947.5

# Row 132 (WVR) - leve item G132
$ws.Range("H132").Value = 3088.1956
$ws.Range("I132").Value = 625.62067
$ws.Range("J132").Value = 7289.0586
$ws.Range("K132").Value = 1876.86201
$ws.Range("L132").Value = 21867.1758
$ws.Range("M132").Value = 653.1379899999999
$ws.Range("N132").Value = -26927.1758
